# Chop.Calc.xlsx edit: rename "Index" column to "i", renumber the index
# column to be zero-based, and narrow column A to fit the new header text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header of column A (also updates the ListObject/table column name
# and the shared-string table automatically).
$ws.Range("A1").Value = "i"

# Renumber the index column A2:A503 from 1..502 down to 0..501 (zero-based).
$n = 502
$arr = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) {
    $arr[$i,0] = $i
}
$ws.Range("A2:A503").Value = $arr

# Narrow column A (was sized for "Index", now fits the shorter "i" header).
$ws.Columns(1).ColumnWidth = 3.1666666666666665
